# Apply "contingencies with rene fine" edit:
# - Insert two new line entries ("line7", "line8") right after "line6" in the
#   name sequence (B column), which pushes the former extr1..extr8 rows down
#   by two positions (extr1..extr6 become rows 10-15, extr7/extr8 become new
#   rows 16-17), and updates the from_bus/to_bus/in_service values for rows
#   8-17 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state for rows 2..17 (A index, name, from_bus, to_bus, in_service)
$rows = @(
    @(0,  "line1", 7,  9,  $true),
    @(1,  "line2", 9,  8,  $true),
    @(2,  "line3", 8,  10, $true),
    @(3,  "line4", 8,  11, $true),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $false),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $true),
    @(11, "extr4", 7,  8,  $true),
    @(12, "extr5", 9,  11, $false),
    @(13, "extr6", 7,  11, $false),
    @(14, "extr7", 5,  7,  $false),
    @(15, "extr8", 8,  5,  $false)
)

$sourceStyleCell = $ws.Cells.Item(2, 1)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    if ($r -gt 15) {
        $sourceStyleCell.Copy()
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
